# 📊 Horarios actualizados Línea 141 - 1104
# Update the three schedule sheets with the newly scraped data
# (new "last updated" timestamp, refreshed row counts and row contents,
# plus a new arrival row that appears on sheet 1 and sheet 2).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:49:51"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = "02:49:51"
$ws1.Range("B6").Value = "02:57"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 8
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = "02:49:51"
$ws1.Range("B7").Value = "04:01"
$ws1.Range("C7").Value = "81_EL PELIGRO"
$ws1.Range("D7").Value = 72
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = "02:49:51"
$ws1.Range("B8").Value = "04:03"
$ws1.Range("C8").Value = "14_ABASTO"
$ws1.Range("D8").Value = 74
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = "02:49:51"
$ws1.Range("B9").Value = "04:45"
$ws1.Range("C9").Value = "215A_EL PATO"
$ws1.Range("D9").Value = 116
$ws1.Range("E9").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:49:51"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A6").Value = "02:49:51"
$ws2.Range("B6").Value = "02:57"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 8
$ws2.Range("E6").Value = "LP1912"

$ws2.Range("A7").Value = "02:49:51"
$ws2.Range("B7").Value = "04:45"
$ws2.Range("C7").Value = "215A_EL PATO"
$ws2.Range("D7").Value = 116
$ws2.Range("E7").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:49:51"
